$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15-19 down to 16-20.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with data (copy constant columns from row 16, the row
# that used to be row 15 before the shift, and set the new/changed values).
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = 44673
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 100112026
$ws.Range("G15").Value = "Haba"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 900
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1400
$ws.Range("M15").Value = 1350
$ws.Range("N15").Value = "$/kilo"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 1350
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"
